$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(48, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(48, 2).Value = 'Mario Leo'
$ws.Cells.Item(48, 3).Value = 'venta'
$ws.Cells.Item(48, 4).Value = 'Lactomayma 22% Peletizada'
$ws.Cells.Item(48, 5).Value = '1 productos'
$ws.Cells.Item(48, 6).Value = 'Q. 198.0'
$ws.Cells.Item(48, 7).Value = 'Q. 198.0'
$ws.Cells.Item(48, 8).Value = 'Q. 3.0'
$ws.Cells.Item(48, 9).Value = 'Q. 184.0'
$ws.Cells.Item(48, 10).Value = 'Q. 11.0'

$ws.Cells.Item(49, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(49, 2).Value = 'Mario Leo'
$ws.Cells.Item(49, 3).Value = 'credito'
$ws.Cells.Item(49, 4).Value = 'Lactomayma'
$ws.Cells.Item(49, 5).Value = '7 productos'
$ws.Cells.Item(49, 6).Value = 'Q. 195.0'
$ws.Cells.Item(49, 7).Value = 'Q. 1365.0'
$ws.Cells.Item(49, 8).Value = 'Q. 21.0'
$ws.Cells.Item(49, 9).Value = 'Q. 1253.0'
$ws.Cells.Item(49, 10).Value = 'Q. 91.0'

$ws.Cells.Item(50, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(50, 2).Value = 'Clientes Varios'
$ws.Cells.Item(50, 3).Value = 'venta'
$ws.Cells.Item(50, 4).Value = 'Pollo Inicio Pelatizado'
$ws.Cells.Item(50, 5).Value = '30 productos'
$ws.Cells.Item(50, 6).Value = 'Q. 255.0'
$ws.Cells.Item(50, 7).Value = 'Q. 7650.0'
$ws.Cells.Item(50, 8).Value = 'Q. 90.0'
$ws.Cells.Item(50, 9).Value = 'Q. 7200.0'
$ws.Cells.Item(50, 10).Value = 'Q. 360.0'

$ws.Cells.Item(51, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(51, 2).Value = 'Leonardo'
$ws.Cells.Item(51, 3).Value = 'credito'
$ws.Cells.Item(51, 4).Value = 'Ponedora fase 2'
$ws.Cells.Item(51, 5).Value = '8 productos'
$ws.Cells.Item(51, 6).Value = 'Q. 234.0'
$ws.Cells.Item(51, 7).Value = 'Q. 1872.0'
$ws.Cells.Item(51, 8).Value = 'Q. 24.0'
$ws.Cells.Item(51, 9).Value = 'Q. 1840.0'
$ws.Cells.Item(51, 10).Value = 'Q. 8.0'

$ws.Cells.Item(52, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(52, 2).Value = 'Clientes Varios'
$ws.Cells.Item(52, 3).Value = 'venta'
$ws.Cells.Item(52, 4).Value = 'Super Vaca lechera'
$ws.Cells.Item(52, 5).Value = '10 productos'
$ws.Cells.Item(52, 6).Value = 'Q. 202.0'
$ws.Cells.Item(52, 7).Value = 'Q. 2020.0'
$ws.Cells.Item(52, 8).Value = 'Q. 30.0'
$ws.Cells.Item(52, 9).Value = 'Q. 1980.0'
$ws.Cells.Item(52, 10).Value = 'Q. 10.0'

$ws.Cells.Item(53, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(53, 2).Value = 'Clientes Varios'
$ws.Cells.Item(53, 3).Value = 'venta'
$ws.Cells.Item(53, 4).Value = 'Super Vaca lechera'
$ws.Cells.Item(53, 5).Value = '1 productos'
$ws.Cells.Item(53, 6).Value = 'Q. 202.0'
$ws.Cells.Item(53, 7).Value = 'Q. 202.0'
$ws.Cells.Item(53, 8).Value = 'Q. 3.0'
$ws.Cells.Item(53, 9).Value = 'Q. 198.0'
$ws.Cells.Item(53, 10).Value = 'Q. 1.0'

$ws.Cells.Item(54, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(54, 2).Value = 'Clientes Varios'
$ws.Cells.Item(54, 3).Value = 'venta'
$ws.Cells.Item(54, 4).Value = 'Pollo Inicio Pelatizado'
$ws.Cells.Item(54, 5).Value = '100 productos'
$ws.Cells.Item(54, 6).Value = 'Q. 255.0'
$ws.Cells.Item(54, 7).Value = 'Q. 25500.0'
$ws.Cells.Item(54, 8).Value = 'Q. 300.0'
$ws.Cells.Item(54, 9).Value = 'Q. 24000.0'
$ws.Cells.Item(54, 10).Value = 'Q. 1200.0'

$ws.Cells.Item(55, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(55, 2).Value = 'Clientes Varios'
$ws.Cells.Item(55, 3).Value = 'venta'
$ws.Cells.Item(55, 4).Value = 'Lactomayma 22% Peletizada'
$ws.Cells.Item(55, 5).Value = '100 productos'
$ws.Cells.Item(55, 6).Value = 'Q. 190.0'
$ws.Cells.Item(55, 7).Value = 'Q. 19000.0'
$ws.Cells.Item(55, 8).Value = 'Q. 300.0'
$ws.Cells.Item(55, 9).Value = 'Q. 18400.0'
$ws.Cells.Item(55, 10).Value = 'Q. 300.0'

$ws.Cells.Item(56, 1).Value = '27 de octubre del 2024'
$ws.Cells.Item(56, 2).Value = 'Resumen'
$ws.Cells.Item(56, 3).Value = 'total'
$ws.Cells.Item(56, 4).Value = 'del'
$ws.Cells.Item(56, 5).Value = 'dia'
$ws.Cells.Item(56, 6).Value = '-'
$ws.Cells.Item(56, 7).Value = 'Q. 54570.0'
$ws.Cells.Item(56, 8).Value = 'Q. 726.0'
$ws.Cells.Item(56, 9).Value = 'Q. 51962.0'
$ws.Cells.Item(56, 10).Value = 'Q. 1882.0'
